$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1834676
$ws.Range("I111").Value = 2200811.5
$ws.Range("J111").Value = 3999
$ws.Range("K111").Value = 6602434.5
$ws.Range("L111").Value = 11997
$ws.Range("M111").Value = -6599367.5
$ws.Range("N111").Value = -18131

$ws.Range("H135").Value = 2204.1765
$ws.Range("I135").Value = 819.4286
$ws.Range("K135").Value = 7374.8574
$ws.Range("M135").Value = -4839.8574

$ws.Range("H137").Value = 4147.533
$ws.Range("I137").Value = 2414.5518
$ws.Range("J137").Value = 7288.5625
$ws.Range("K137").Value = 7243.655400000001
$ws.Range("L137").Value = 21865.6875
$ws.Range("M137").Value = -4693.655400000001
$ws.Range("N137").Value = -26965.6875

$ws.Range("H138").Value = 5380.353
$ws.Range("I138").Value = 3975.923
$ws.Range("J138").Value = 5860.816
$ws.Range("K138").Value = 11927.769
$ws.Range("L138").Value = 17582.448
$ws.Range("M138").Value = -6787.769
$ws.Range("N138").Value = -27862.448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 816
$ws.Range("I25").Value = 816
$ws.Range("K25").Value = 816
$ws.Range("M25").Value = -414

$ws.Range("H32").Value = 2840.375
$ws.Range("I32").Value = 2476.209
$ws.Range("K32").Value = 2476.209
$ws.Range("M32").Value = -2189.209

$ws.Range("H61").Value = 4729.6924
$ws.Range("I61").Value = 3185.875
$ws.Range("K61").Value = 3185.875
$ws.Range("M61").Value = -2973.875

$ws.Range("H92").Value = 12532000
$ws.Range("J92").Value = 42666.668
$ws.Range("L92").Value = 42666.668
$ws.Range("N92").Value = -47658.668

$ws.Range("H102").Value = 3141
$ws.Range("I102").Value = 2986.0833
$ws.Range("K102").Value = 2986.0833
$ws.Range("M102").Value = -1364.0833

$ws.Range("H122").Value = 4215.1924
$ws.Range("I122").Value = 2867.5293
$ws.Range("K122").Value = 8602.5879
$ws.Range("M122").Value = -6152.5879

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws.Range("H135").Value = 89537
$ws.Range("J135").Value = 89537
$ws.Range("L135").Value = 89537
$ws.Range("N135").Value = -99677

$ws.Range("H136").Value = 4729.6924
$ws.Range("I136").Value = 3185.875
$ws.Range("K136").Value = 9557.625
$ws.Range("M136").Value = -7007.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3307.9473
$ws.Range("I99").Value = 3170.7646
$ws.Range("K99").Value = 3170.7646
$ws.Range("M99").Value = -1672.7646

$ws.Range("H105").Value = 3761.818
$ws.Range("I105").Value = 1371.909
$ws.Range("J105").Value = 4956.773
$ws.Range("K105").Value = 1371.909
$ws.Range("L105").Value = 4956.773
$ws.Range("M105").Value = 375.0909999999999
$ws.Range("N105").Value = -8450.773000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6685.3335
$ws.Range("I31").Value = 2112
$ws.Range("K31").Value = 2112
$ws.Range("M31").Value = -1817

$ws.Range("H34").Value = 6685.3335
$ws.Range("I34").Value = 2112
$ws.Range("K34").Value = 2112
$ws.Range("M34").Value = -1910

$ws.Range("H57").Value = 128055
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32372

$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -101856

$ws.Range("H93").Value = 39500
$ws.Range("I93").Value = 39500
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 39500
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -37628
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 6879.75
$ws.Range("I133").Value = 8339.666999999999
$ws.Range("K133").Value = 25019.001
$ws.Range("M133").Value = -19959.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 26819.785
$ws.Range("J15").Value = 26819.785
$ws.Range("L15").Value = 26819.785
$ws.Range("N15").Value = -27395.785

$ws.Range("H51").Value = 60000
$ws.Range("J51").Value = 60000
$ws.Range("L51").Value = 60000
$ws.Range("N51").Value = -61018

$ws.Range("H81").Value = 26819.785
$ws.Range("J81").Value = 26819.785
$ws.Range("L81").Value = 26819.785
$ws.Range("N81").Value = -28815.785

$ws.Range("H84").Value = 26819.785
$ws.Range("J84").Value = 26819.785
$ws.Range("L84").Value = 80459.355
$ws.Range("N84").Value = -90443.355

$ws.Range("H122").Value = 5614.8936
$ws.Range("I122").Value = 5520.1
$ws.Range("J122").Value = 5782.1763
$ws.Range("K122").Value = 16560.3
$ws.Range("L122").Value = 17346.5289
$ws.Range("M122").Value = -14110.3
$ws.Range("N122").Value = -22246.5289

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 373829.62
$ws.Range("I132").Value = 403376.4
$ws.Range("J132").Value = 4495
$ws.Range("K132").Value = 1210129.2
$ws.Range("L132").Value = 13485
$ws.Range("M132").Value = -1207599.2
$ws.Range("N132").Value = -18545

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 28750
$ws.Range("J97").Value = 28750
$ws.Range("L97").Value = 28750
$ws.Range("N97").Value = -30732

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H136").Value = 4596.1333
$ws.Range("I136").Value = 3370.375
$ws.Range("K136").Value = 10111.125
$ws.Range("M136").Value = -7561.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 41999.5
$ws.Range("J70").Value = 40000
$ws.Range("L70").Value = 40000
$ws.Range("N70").Value = -40630

$ws.Range("H73").Value = 41999.5
$ws.Range("J73").Value = 40000
$ws.Range("L73").Value = 40000
$ws.Range("N73").Value = -42184

$ws.Range("H126").Value = 5565.077
$ws.Range("I126").Value = 6008
$ws.Range("J126").Value = 5185.4287
$ws.Range("K126").Value = 18024
$ws.Range("L126").Value = 15556.2861
$ws.Range("M126").Value = -15554
$ws.Range("N126").Value = -20496.2861
